$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.760.63"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.295.34"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "2.653.06"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "2.281.20"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "42.672.61"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.46%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").Value = "1.999.28"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "2.520.89"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
